# CheckCredibility returns a dict object
#
# The underlying Python helper (CheckCredibility) used to return values
# from an ordered structure; now it returns a dict, which changed the
# order that rows got written into the "COMPARING" sheet. Row 4's data
# moved up to row 2, and the old rows 2 and 3 shifted down to 3 and 4
# (a 3-row rotation). This script reproduces the resulting workbook
# state: the rotated COMPARING rows, the refreshed column widths for
# the newly-visible J:L block, and the active-sheet/selection churn
# that came along with the author's last save (COMPARING became the
# active tab instead of the SPACY sheet).

$wb = $excel.ActiveWorkbook

$wsComparing = $wb.Worksheets.Item("COMPARING")
$wsSpacy = $wb.Worksheets.Item("SPACY Names Entity Recognition")

# --- Rotate the data rows 2,3,4 on COMPARING -------------------------------
# new row2 <- old row4, new row3 <- old row2, new row4 <- old row3
# (value + number format + fill all travel together, so we physically
# copy the cells rather than only copying .Value)
$cols = "B","D","E","F","G","H","J","K","L"

foreach ($col in $cols) {
  $wsComparing.Range($col + "4").Copy($wsComparing.Range($col + "200"))
}
foreach ($col in $cols) {
  $wsComparing.Range($col + "2").Copy($wsComparing.Range($col + "201"))
}
foreach ($col in $cols) {
  $wsComparing.Range($col + "3").Copy($wsComparing.Range($col + "202"))
}

foreach ($col in $cols) {
  $wsComparing.Range($col + "200").Copy($wsComparing.Range($col + "2"))
}
foreach ($col in $cols) {
  $wsComparing.Range($col + "201").Copy($wsComparing.Range($col + "3"))
}
foreach ($col in $cols) {
  $wsComparing.Range($col + "202").Copy($wsComparing.Range($col + "4"))
}

# Row 2 no longer has a subjectivity label (that moved out with old row4,
# which had nothing in F) - clear the leftover copy explicitly.
$wsComparing.Range("F2").Clear()

# Drop the staging rows used for the rotation.
$wsComparing.Range("B200:L202").Clear()

# --- Column widths for I:L (now visible/relevant again) -------------------
# Matches the widths Excel's own AutoFit produced for the author (computed
# from the "ColumnWidth input" <-> "stored width" offset this host uses).
$wsComparing.Columns("I:I").ColumnWidth = 3.7213541666666665
$wsComparing.Columns("J:J").ColumnWidth = 6.166666666666667
$wsComparing.Columns("K:K").ColumnWidth = 7.166666666666667
$wsComparing.Columns("L:L").ColumnWidth = 6.166666666666667

# --- Active sheet / selection churn ----------------------------------------
# The SPACY sheet was selected before; the author ended up on COMPARING.
$wsSpacy.Activate() | Out-Null
$wsSpacy.Range("J13").Select() | Out-Null

$wsComparing.Activate() | Out-Null
$wsComparing.Range("O17").Select() | Out-Null
